# The deck's "datetimeFigureOut" fields (auto date placeholders) were
# cached from 25/04/2018; the file was touched again a day later, so
# every cached date field (slide master, every slide layout, and the
# notes master) needs to roll forward from the 25th to the 26th.
#
# There are no slide-level date fields in this deck (only inherited
# ones on the master/layouts/notes master), so only those containers
# are touched.

$p = $ppt.ActivePresentation

function Update-DateField($shape) {
    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    $t = $tr.Text
    if ($t -eq $null) { return }
    if ($t -match "^25/04/2018$") {
        $tr.Text = "26/04/2018"
    } elseif ($t -match "^4/25/2018$") {
        $tr.Text = "4/26/2018"
    }
}

function Update-Container($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        Update-DateField $container.Shapes.Item($i)
    }
}

# Slide master
Update-Container $p.SlideMaster

# Every slide layout inherits its own cached copy of the field
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-Container $layouts.Item($j)
}

# Notes master (en-GB locale -> "DD/MM/YYYY" formatted date)
Update-Container $p.NotesMaster
